$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some 'Price' values are plain decimal numbers (e.g. '536.17').
# Assigning such a string straight to Range.Value lets Excel's
# automatic type inference turn it into a real Number, but the
# source data stores every cell in this sheet as text (note values
# like '59.496.40' or '2.612.52' that are NOT valid numbers at all).
# To keep these cells as text - matching the rest of the sheet -
# we stage the value in a helper cell that's explicitly formatted
# as Text, then use Copy/PasteSpecial (values only) to move just
# the literal text into the target cell without touching its own
# number format/style.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range('D2').Value = '59.496.40'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.638.71'
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('E4').Value = '  +0.30%  '
$helper.Value = '536.17'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range('E5').Value = '  -0.32%  '
$helper.Value = '144.96'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range('E6').Value = '  +2.84%  '
$helper.Value = '0.999'
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range('E7').Value = '  +0.00%  '
$helper.Value = '0.570'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = '2.652.47'
$ws.Range('E9').Value = '  +1.36%  '
$helper.Value = '6.64'
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range('E10').Value = '  +2.62%  '
$helper.Value = '0.102'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range('E11').Value = '  -1.06%  '
$helper.Value = '0.337'
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '3.110.70'
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').Value = '59.419.42'
$ws.Range('E15').Value = '  +0.09%  '
$helper.Value = '21.09'
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$helper.Value = '0.0000134'
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.612.52'
$ws.Range('E18').Value = '  +0.27%  '
$helper.Value = '339.89'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('E20').Value = '  +0.75%  '
$helper.Value = '10.35'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('E22').Value = '  -1.58%  '
$ws.Range('E23').Value = '  +0.14%  '
$helper.Value = '67.08'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('E26').Value = '  -1.45%  '
$helper.Value = '0.998'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range('E27').Value = '  -0.14%  '
$helper.Value = '7.27'
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').Value = '0.0₃0745'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  +0.33%  '
$helper.Value = '5.83'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range('E32').Value = '  -0.43%  '
$helper.Value = '18.90'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range('E33').Value = '  +0.11%  '
$helper.Value = '151.51'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range('E34').Value = '  +1.78%  '
$helper.Value = '3.99'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('E37').Value = '  -0.46%  '
$helper.Value = '0.833'
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('E39').Value = '  -0.98%  '
$helper.Value = '288.54'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range('E40').Value = '  +4.16%  '
$ws.Range('E41').Value = '  +0.90%  '
$helper.Value = '0.999'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  +3.42%  '
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').Value = '1.972.52'
$ws.Range('E48').Value = '  +1.06%  '
$helper.Value = '0.0225'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range('E49').Value = '  +0.81%  '
$helper.Value = '4.54'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range('E50').Value = '  +0.42%  '
$helper.Value = '18.26'
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range('E51').Value = '  -0.62%  '

$helper.Clear()
$excel.CutCopyMode = $false
